# display_template.xlsx update
# - fix a few display-name typos/wording
# - add a new "Morrisons FTG Kit (FTG)" row at the bottom of the table
# - extend the table / autofilter / filter-database named ranges to the new last row
# - minor cosmetic row-height / column-width / selection tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Small text corrections to existing rows (shared-string edits in the diff)
# ---------------------------------------------------------------------------
$ws.Range("A47").Value  = "o. Box Display"
$ws.Range("A92").Value  = "HO Agreed Gondola End"
$ws.Range("A93").Value  = "Non HO Gondola End"
$ws.Range("A95").Value  = "Non HO Agreed shelf display"

# ---------------------------------------------------------------------------
# 2. Append the new row (105) with the new "Morrisons FTG Kit (FTG)" entry
# ---------------------------------------------------------------------------
$ws.Range("A105").Value = "Morrisons FTG Kit (FTG)"
$ws.Range("B105").Value = "Shelf"
$ws.Range("C105").Value = "N/A"
$ws.Range("D105").Value = "No"
$ws.Range("E105").Value = "Yes"
$ws.Range("F105").Value = "No"

# match the formatting used by the rest of the table body (copy down from the
# row above rather than hand-building styles so existing style records are
# reused the same way Excel/LO would when you fill a new row in)
$ws.Range("A2").Copy()
$ws.Range("A105").PasteSpecial(-4122)
$ws.Range("A105").Style = "Normal"

$ws.Range("B2:F2").Copy()
$ws.Range("B105:F105").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Re-apply the AutoFilter over the grown range (A1:F105) and refresh the
#    worksheet dimension so the sheet picks up the new row.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:F105").AutoFilter()

# ---------------------------------------------------------------------------
# 4. Fix up the workbook-level defined names that track the filter database
#    (these are not auto-synced by AutoFilter() in this host, unlike real
#    Excel, so update them by hand to mirror what Excel/LO would do).
# ---------------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$F`$105"
$wb.Names.Item(4).Delete()
$wb.Names.Item(2).Delete()
$newFilterName = $ws.Names.Add("TempFilterDB", "=Sheet1!`$A`$1:`$F`$105")
$newFilterName.Name = "_xlnm._FilterDatabase"

# ---------------------------------------------------------------------------
# 5. Cosmetic tweaks: row heights, column widths and the saved selection.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 39
$ws.Rows.Item(11).RowHeight = 12.8

$ws.Columns.Item(1).ColumnWidth = 30.89
$ws.Columns.Item(3).ColumnWidth = 14.42
$ws.Columns.Item(5).ColumnWidth = 11.45

$ws.Activate()
$ws.Range("C92").Select()
